# "Generate Report for Handoff" — mark the c0848784 entry as ready for
# handoff (it used to be "Handed back: in sync with en-US") and record the
# new handoff timestamps + an error detail message explaining that the
# handback file that was picked up is stale.

$wb = $excel.ActiveWorkbook

$status = "Ready for handoff"
$errorDetail = "The version of handback file is not the latest, current: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/eb450ffc03a749b33404fdb9a4568368387fc64b/e2e/c0848784-12d9-43de-9a7a-2d97d0b82dd1.md, latest: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/a44cabead72b817eda9dcd29ecd45c4e659a255b/e2e/c0848784-12d9-43de-9a7a-2d97d0b82dd1.md."

# ---- Overview sheet: row 3 is the c0848784 file ----
$overview = $wb.Worksheets.Item("Overview")
$overview.Range("E3").Value = $status
$overview.Range("F3").Value = $status
$overview.Range("G3").Value = "2016-09-05 23:01:22"

# ---- zh-cn sheet: row 3 is the c0848784 file ----
$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Range("C3").Value = $status
$zhcn.Range("H3").Value = "2016-09-05 23:01:18"
$zhcn.Range("P3").Value = $errorDetail
$zhcn.Columns.Item(16).ColumnWidth = 39.17

# ---- de-de sheet: row 3 is the c0848784 file ----
$dede = $wb.Worksheets.Item("de-de")
$dede.Range("C3").Value = $status
$dede.Range("H3").Value = "2016-09-05 23:01:22"
$dede.Range("P3").Value = $errorDetail
$dede.Columns.Item(16).ColumnWidth = 39.17
